$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.396.43'
$ws.Range("E2").Value = '  +8.28%  '
$ws.Range("D3").Value = '1.675.52'
$ws.Range("E3").Value = '  +4.07%  '
$ws.Range("D4").Value = '''1.007'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").Value = '''0.9999'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.93%  '
$ws.Range("D6").Value = '''305.93'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '''0.3697'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("D8").Value = '''0.3418'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.45%  '
$ws.Range("D9").Value = '''47.77'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +12.92%  '
$ws.Range("D10").Value = '''1.156'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.61%  '
$ws.Range("D11").Value = '''0.07221'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.30%  '
$ws.Range("D12").Value = '''1.003'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("D13").Value = '''6.077'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.81%  '
$ws.Range("D14").Value = '''20.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.47%  '
$ws.Range("D15").Value = '''6.707'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.92%  '
$ws.Range("D16").Value = '1.676.33'
$ws.Range("E16").Value = '  +4.15%  '
$ws.Range("D17").Value = '''0.00001101'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("D18").Value = '''1.000'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.90%  '
$ws.Range("D19").Value = '''0.06660'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("D20").Value = '''80.65'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.29%  '
$ws.Range("D21").Value = '''16.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.72%  '
$ws.Range("D22").Value = '''6.069'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.03%  '
$ws.Range("D23").Value = '''12.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.81%  '
$ws.Range("D24").Value = '24.348.94'
$ws.Range("E24").Value = '  +8.00%  '
$ws.Range("D25").Value = '''2.422'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.99%  '
$ws.Range("E26").Value = '  -13.67%  '
$ws.Range("D27").Value = '''2.635'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.81%  '
$ws.Range("D28").Value = '''152.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.88%  '
$ws.Range("D29").Value = '''19.32'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.09%  '
$ws.Range("D30").Value = '1.861.24'
$ws.Range("E30").Value = '  +3.86%  '
$ws.Range("D31").Value = '''127.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.94%  '
$ws.Range("D32").Value = '''6.251'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.11%  '
$ws.Range("D33").Value = '''4.060'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.27%  '
$ws.Range("D34").Value = '''0.9631'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.39%  '
$ws.Range("D35").Value = '''0.08406'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.04%  '
$ws.Range("E36").Value = '  +1.16%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '''0.06423'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.88%  '
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").Value = '''12.21'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.41%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '''5.293'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.79%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '''8.793'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.24%  '
$ws.Range("D41").Value = '''0.02301'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.67%  '
$ws.Range("D42").Value = '''1.234'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.81%  '
$ws.Range("D43").Value = '''0.2083'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.04%  '
$ws.Range("D44").Value = '''0.6091'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.92%  '
$ws.Range("D45").Value = '''0.9999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.94%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '''3.757'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.16%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''13.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.28%  '
$ws.Range("D48").Value = '''0.5878'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.24%  '
$ws.Range("D49").Value = '''126.63'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("D50").Value = '''2.005'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.76%  '
$ws.Range("D51").Value = '''0.07154'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.87%  '
